$wb = $excel.ActiveWorkbook

# "Generate Report for handback": refresh the handoff / handback timestamps
# for the file that was just handed back (5e39bc3a-...) in each language sheet.
# The other file's row (a7de36fb-...) keeps its previous timestamps.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-18 03:46:59"
$wsZhCn.Range("G2").Value = "2016-01-18 03:47:44"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-18 03:47:11"
$wsDeDe.Range("G2").Value = "2016-01-18 03:48:01"
